# Add two new columns "I0" (I) and "IF" (J) with per-row numeric data,
# matching the header style already used by the existing header row (B1:H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: copy the existing header formatting (bold, bordered,
# centered) from H1 onto I1:J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row data for columns I ("I0") and J ("IF"), rows 2-74.
$iVals = @(6,5,7,8,9,7,9,7,5,8,8,9,8,9,7,3,7,8,7,6,7,8,8,7,8,7,8,8,8,7,8,7,8,10,8,7,8,7,7,9,9,9,9,9,9,7,8,8,9,7,8,8,7,7,8,7,8,8,8,8,7,7,8,9,9,10,7,8,9,8,8,7,4)
$jVals = @(6,7,8,9,9,7,9,7,5,8,8,9,8,9,7,4,7,8,7,6,7,8,8,7,8,7,8,8,8,7,8,7,8,10,8,7,8,7,7,9,9,9,9,9,9,7,8,8,9,7,8,8,7,7,8,8,8,8,8,8,7,7,9,9,9,10,7,8,9,8,8,7,4)

$firstRow = 2
for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $r = $firstRow + $idx
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
